$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the analytical-solution formula in AV2 and propagate as a shared
# formula across AV3:AV18 (matches original layout: AV2 standalone,
# AV3 anchors the shared formula used by AV4:AV18).
$ws.Range("AV2").Formula = "=(AQ2*X2-4*X2*(1-0.01*P2-2*0.01*AF2)*(0.08/0.4)/(-0.08/0.4*0.01*P2-(2*0.08/0.4+3)*0.01*AF2+0.08/0.4+1))/1000"
$ws.Range("AV3:AV18").Formula = "=(AQ3*X3-4*X3*(1-0.01*P3-2*0.01*AF3)*(0.08/0.4)/(-0.08/0.4*0.01*P3-(2*0.08/0.4+3)*0.01*AF3+0.08/0.4+1))/1000"

# Update the selection state to match the author's UI action: select the
# full column AV with active cell AV1.
$ws.Range("AV1:AV1048576").Select()
